$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A and append the new ticker rows after it
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Cells.Item($lastRow + 1, 1).Value = "IMX-USD"
$ws.Cells.Item($lastRow + 2, 1).Value = "GRT-USD"
